$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-10-07 Tuesday"; new = "2025-10-08 Wednesday"},
    @{old = "57÷9="; new = "22÷5="},
    @{old = "51÷2="; new = "96÷6="},
    @{old = "80÷6="; new = "50÷4="},
    @{old = "36÷7="; new = "99÷9="},
    @{old = "34÷5="; new = "50÷8="},
    @{old = "21÷3="; new = "26÷5="},
    @{old = "57÷8="; new = "81÷7="},
    @{old = "86÷5="; new = "75÷5="},
    @{old = "85÷9="; new = "35÷7="},
    @{old = "51÷9="; new = "13÷6="},
    @{old = "58÷3="; new = "66÷4="},
    @{old = "23÷9="; new = "73÷8="},
    @{old = "99÷8="; new = "54÷6="},
    @{old = "78÷4="; new = "15÷9="},
    @{old = "76÷7="; new = "84÷4="},
    @{old = "88÷6="; new = "17÷8="},
    @{old = "75÷4="; new = "15÷3="},
    @{old = "52÷4="; new = "94÷4="},
    @{old = "22÷3="; new = "54÷6="},
    @{old = "90÷8="; new = "16÷4="},
    @{old = "21÷6="; new = "37÷4="},
    @{old = "43÷5="; new = "66÷4="},
    @{old = "75÷9="; new = "41÷3="},
    @{old = "16÷3="; new = "10÷9="},
    @{old = "74÷2="; new = "48÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2) | Out-Null
}
